$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Overview sheet: Status moved from "Ready for handoff" to
# "Handed back: in sync with en-US" for both the zh-cn and de-de columns.
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"

# Widen the two status columns to fit the longer text.
$overview.Columns("E").ColumnWidth = 29.1
$overview.Columns("F").ColumnWidth = 29.1

# ---------------------------------------------------------------------------
# zh-cn sheet: record the handback target file, handback file and handback
# datetime now that the round trip has completed.
# ---------------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("C2").Value = "Handed back: in sync with en-US"

$zhcnTarget = $zhcn.Range("I2")
$zhcnTarget.Value = "e3a8775b-65f1-4368-91a2-95eef1b2486e.md"
$zhcn.Hyperlinks.Add($zhcnTarget, "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e8bcc410dce47b6517a128402289c4573ac15f44/e2e/e3a8775b-65f1-4368-91a2-95eef1b2486e.md", "", "", "e3a8775b-65f1-4368-91a2-95eef1b2486e.md")

$zhcn.Range("J2").Value = "e3a8775b-65f1-4368-91a2-95eef1b2486e.f50278a990b0bd43566c460405169ea95900c7d9.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-09-01 21:08:19"

$zhcn.Columns("C").ColumnWidth = 29.1
$zhcn.Columns("I").ColumnWidth = 39.2
$zhcn.Columns("J").ColumnWidth = 39.2

# ---------------------------------------------------------------------------
# de-de sheet: same bookkeeping as zh-cn, but with the de-de handback file
# and its own handback timestamp.
# ---------------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("C2").Value = "Handed back: in sync with en-US"

$dedeTarget = $dede.Range("I2")
$dedeTarget.Value = "e3a8775b-65f1-4368-91a2-95eef1b2486e.md"
$dede.Hyperlinks.Add($dedeTarget, "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e8bcc410dce47b6517a128402289c4573ac15f44/e2e/e3a8775b-65f1-4368-91a2-95eef1b2486e.md", "", "", "e3a8775b-65f1-4368-91a2-95eef1b2486e.md")

$dede.Range("J2").Value = "e3a8775b-65f1-4368-91a2-95eef1b2486e.f50278a990b0bd43566c460405169ea95900c7d9.de-de.xlf"
$dede.Range("K2").Value = "2016-09-01 21:08:26"

$dede.Columns("C").ColumnWidth = 29.1
$dede.Columns("I").ColumnWidth = 39.2
$dede.Columns("J").ColumnWidth = 39.2
